# Avoid interactive "are you sure" prompts (e.g. when deleting a sheet)
$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the "Desarquivamentos Pendentes" sheet entirely
$wb.Worksheets("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
